# Insert a new row at row 77, shifting existing rows 77-179 down to 78-180
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("77:77").Insert()

# Populate the newly inserted row 77 with the new record's data
$ws.Range("A77").Value = 10
$ws.Range("B77").Value = "Vega Modelo de Temuco"
$ws.Range("C77").Value = "La Araucanía"
$ws.Range("D77").Value = 44413
$ws.Range("E77").Value = 9
$ws.Range("F77").Value = 100112008
$ws.Range("G77").Value = "Coliflor"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 1000
$ws.Range("L77").Value = 1000
$ws.Range("M77").Value = 1000
$ws.Range("N77").Value = "$/unidad"
$ws.Range("O77").Value = "Región Metropolitana"
$ws.Range("P77").Value = 1000
$ws.Range("Q77").Value = 1
$ws.Range("R77").Value = "Hortaliza"

# Ensure the D77 cell keeps the date/datetime number format used by the rest of column D
$ws.Range("D77").NumberFormat = $ws.Range("D78").NumberFormat
